$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 3-5: normalize the BẮT ĐẦU / KẾT THÚC date strings to 20/09/23 / 26/09/23 ---
$ws.Range("C3").Value = "20/09/23"
$ws.Range("D3").Value = "26/09/23"

$ws.Range("C4").Value = "20/09/23"
$ws.Range("D4").Value = "26/09/23"
$ws.Range("F4").Value = "Hoàn thành"

$ws.Range("C5").Value = "20/09/23"
$ws.Range("D5").Value = "26/09/23"

# --- Rows 8-16: fill in BẮT ĐẦU / KẾT THÚC / THÀNH VIÊN / TÌNH TRẠNG ---
# Seed the shared-string table so "Nguyễn Phạm Nhật Minh" is registered before
# "Bùi Phi Long" (matches the row-9-before-row-8 first-use order of the target file).
$ws.Range("E9").Value = "Nguyễn Phạm Nhật Minh"
$ws.Range("E8").Value = "Bùi Phi Long"

$members = @{
    8  = "Bùi Phi Long"
    9  = "Nguyễn Phạm Nhật Minh"
    10 = "Nguyễn Phạm Nhật Minh"
    11 = "Bùi Phi Long"
    12 = "Nguyễn Phạm Nhật Minh"
    13 = "Bùi Phi Long"
    14 = "Bùi Phi Long"
    15 = "Nguyễn Phạm Nhật Minh"
    16 = "Bùi Phi Long"
}

# Set up the "seed" cells (row 8) with the real formats first, then copy the
# formatting down to rows 9-16 so every row shares the same two style records
# instead of each row minting its own (mirrors how Excel itself would behave
# if the date columns were formatted once and filled down).
$ws.Range("D8").Value = 45026
$ws.Range("D8").NumberFormat = "mm-dd-yy"
$ws.Range("C8").Value = "'26/09/23"
$ws.Range("C8").NumberFormat = "mm-dd-yy"

$ws.Range("D8").Copy()
$ws.Range("D9:D16").PasteSpecial(-4122)
$ws.Range("D9:D16").Value = 45026

$ws.Range("C8").Copy()
$ws.Range("C9:C16").PasteSpecial(-4122)
$ws.Range("C9:C16").Value = "'26/09/23"

foreach ($row in 8..16) {
    $ws.Cells.Item($row, 5).Value = $members[$row]
    $ws.Cells.Item($row, 6).Value = "Hoàn thành"
}

# --- Update the stored selection to match the author's last-saved cursor position ---
$ws.Range("G16").Select()
